$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.214.71'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '2.244.24'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.14'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.21'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.602'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.41'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0959'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.95'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.56%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("D14").Value = '2.580.60'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.39'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.841'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").Value = '2.224.02'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").Value = '42.116.26'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("E19").Value = '  -4.38%  '
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.75'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.29'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +8.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.42'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -6.72%  '
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("E26").Value = '  -2.64%  '
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.28'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.98%  '
$ws.Range("E29").Value = '  -2.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.37'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.64'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.65'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0805'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.80'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.56%  '
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("E36").Value = '  -6.43%  '
$ws.Range("E37").Value = '  -4.82%  '
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("E39").Value = '  -1.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.14'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.94'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.22%  '
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("E44").Value = '  -1.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.31'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("E46").Value = '  -1.73%  '
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.34'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.89%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.17'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("D51").Value = '2.451.73'
$ws.Range("E51").Value = '  +0.05%  '
